# Update "想去人数" (want-to-go count) figures in column F
# for both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 61
    7  = 1253
    8  = 1530
    15 = 105
    19 = 1728
    23 = 666
    26 = 4158
    28 = 267
    29 = 1084
    30 = 485
    32 = 525
    34 = 241
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
